# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# cells on the per-language handback sheets, as produced by re-running the
# "Generate Report for handback" step.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-28 11:31:53"
$wsZhCn.Range("G5").Value = "2016-01-28 11:32:36"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-28 11:32:06"
$wsDeDe.Range("G5").Value = "2016-01-28 11:33:00"
